$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix C35: was stored as text, should be numeric ---
$ws.Range("C35").Value = 21691136843814

# --- Row 36 ---
$ws.Range("A36").Value = "Large Images"
$ws.Range("B36").Value = "segmentation"
$ws.Range("C36").Value = 21691136843814
$ws.Range("D36").Value = "9.30 minutes"
$ws.Range("E36").Value = "SegFormer-[14M]"
$ws.Range("F36").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G36").Value = "'73.19%"
$ws.Range("H36").Value = "'39.49%"
$ws.Range("I36").Value = "'56.72%"
$ws.Range("J36").Value = "2025-03-09 16:09:59"

# --- Row 37 ---
$ws.Range("A37").Value = "Large Images"
$ws.Range("B37").Value = "segmentation"
$ws.Range("C37").Value = 21691136843814
$ws.Range("D37").Value = "9.49 minutes"
$ws.Range("E37").Value = "SegFormer-[14M]"
$ws.Range("F37").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G37").Value = "'72.43%"
$ws.Range("H37").Value = "'45.77%"
$ws.Range("I37").Value = "'56.86%"
$ws.Range("J37").Value = "2025-03-09 16:59:31"

# --- Row 38 ---
$ws.Range("A38").Value = "Large Images"
$ws.Range("B38").Value = "segmentation"
$ws.Range("C38").Value = "'21691136843814"
$ws.Range("D38").Value = "9.41 minutes"
$ws.Range("E38").Value = "SegFormer-[14M]"
$ws.Range("F38").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G38").Value = "'73.27%"
$ws.Range("H38").Value = "'39.25%"
$ws.Range("I38").Value = "'56.69%"
$ws.Range("J38").Value = "2025-03-10 00:45:58"

# --- Row 39 ---
$ws.Range("A39").Value = "Large Images"
$ws.Range("B39").Value = "object-detection"
$ws.Range("C39").Value = "'103610396583974"
$ws.Range("D39").Value = "5.50 minutes"
$ws.Range("E39").Value = "RepPoints-[37M]"
$ws.Range("F39").Value = "{'height': 2000, 'width': 2000, 'paddingValue': 0}"
$ws.Range("G39").Value = "'85.25%"
$ws.Range("H39").Value = "N/A"
$ws.Range("I39").Value = "N/A"
$ws.Range("J39").Value = "2025-03-10 00:59:54"

# --- Row 40 ---
$ws.Range("A40").Value = "Large Images"
$ws.Range("B40").Value = "segmentation"
$ws.Range("C40").Value = "'21691136843814"
$ws.Range("D40").Value = "9.26 minutes"
$ws.Range("E40").Value = "SegFormer-[14M]"
$ws.Range("F40").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G40").Value = "'73.77%"
$ws.Range("H40").Value = "'43.71%"
$ws.Range("I40").Value = "'56.74%"
$ws.Range("J40").Value = "2025-03-10 01:11:24"

# --- Row 41 ---
$ws.Range("A41").Value = "Large Images"
$ws.Range("B41").Value = "segmentation"
$ws.Range("C41").Value = "'21691136843814"
$ws.Range("D41").Value = "10.47 minutes"
$ws.Range("E41").Value = "FastVit-[14M]"
$ws.Range("F41").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G41").Value = "'72.00%"
$ws.Range("H41").Value = "'52.28%"
$ws.Range("I41").Value = "'56.32%"
$ws.Range("J41").Value = "2025-03-10 01:22:11"

# --- Row 42 ---
$ws.Range("A42").Value = "Large Images"
$ws.Range("B42").Value = "segmentation"
$ws.Range("C42").Value = "'1938830"
$ws.Range("D42").Value = "6.44 minutes"
$ws.Range("E42").Value = "FastVit-[14M]"
$ws.Range("F42").Value = "{'height': 2048, 'width': 2048, 'paddingValue': 0}"
$ws.Range("G42").Value = "'19.41%"
$ws.Range("H42").Value = "'6.78%"
$ws.Range("I42").Value = "'24.80%"
$ws.Range("J42").Value = "2025-03-10 01:35:01"

# --- Row 43 ---
$ws.Range("A43").Value = "Large Images"
$ws.Range("B43").Value = "segmentation"
$ws.Range("C43").Value = "'1938830"
$ws.Range("D43").Value = "6.05 minutes"
$ws.Range("E43").Value = "SegFormer-[14M]"
$ws.Range("F43").Value = "{'height': 2048, 'width': 2048, 'paddingValue': 0}"
$ws.Range("G43").Value = "'30.07%"
$ws.Range("H43").Value = "'11.42%"
$ws.Range("I43").Value = "'45.54%"
$ws.Range("J43").Value = "2025-03-10 01:44:12"

# --- Row 44 ---
$ws.Range("A44").Value = "Large Images"
$ws.Range("B44").Value = "object-detection"
$ws.Range("C44").Value = "'102004060440613"
$ws.Range("D44").Value = "11.19 minutes"
$ws.Range("E44").Value = "RtmDet-[9M]"
$ws.Range("F44").Value = "{'height': 3040, 'width': 4056, 'paddingValue': 0}"
$ws.Range("G44").Value = "'70.80%"
$ws.Range("H44").Value = "'68.75%"
$ws.Range("I44").Value = "'66.67%"
$ws.Range("J44").Value = "2025-03-10 01:56:45"

# --- Row 45 ---
$ws.Range("A45").Value = "Large Images"
$ws.Range("B45").Value = "object-detection"
$ws.Range("C45").Value = "'6464689526794"
$ws.Range("D45").Value = "8.38 minutes"
$ws.Range("E45").Value = "RtmDet-[9M]"
$ws.Range("F45").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G45").Value = "'77.68%"
$ws.Range("H45").Value = "'76.92%"
$ws.Range("I45").Value = "'76.92%"
$ws.Range("J45").Value = "2025-03-10 02:08:57"

# --- Row 46 ---
$ws.Range("A46").Value = "Large Images"
$ws.Range("B46").Value = "object-detection"
$ws.Range("C46").Value = "'6464689526794"
$ws.Range("D46").Value = "10.03 minutes"
$ws.Range("E46").Value = "RepPoints-[20M]"
$ws.Range("F46").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G46").Value = "'76.56%"
$ws.Range("H46").Value = "'80.00%"
$ws.Range("I46").Value = "'73.08%"
$ws.Range("J46").Value = "2025-03-10 02:24:36"

# --- Row 47 ---
$ws.Range("A47").Value = "Large Images"
$ws.Range("B47").Value = "object-detection"
$ws.Range("C47").Value = "'6464689526794"
$ws.Range("D47").Value = "12.82 minutes"
$ws.Range("E47").Value = "RepPoints-[37M]"
$ws.Range("F47").Value = "{'height': 6000, 'width': 6000, 'paddingValue': 0}"
$ws.Range("G47").Value = "'73.38%"
$ws.Range("H47").Value = "'71.43%"
$ws.Range("I47").Value = "'78.26%"
$ws.Range("J47").Value = "2025-03-10 02:37:41"
